$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "总计" sheet: insert a new "2022-Q3" row right after the header, shifting
#    the existing quarter rows down by one and renumbering column A.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item(1)

# Push row 4 (2021-Q2) data down into row 5, copying A4's number style along.
$summary.Range("B5").Value = $summary.Range("B4").Value2
$summary.Range("C5").Value = $summary.Range("C4").Value2
$summary.Range("D5").Value = $summary.Range("D4").Value2
$summary.Range("A4").Copy()
$summary.Range("A5").PasteSpecial(-4122)

# Push row 3 (2021-Q3) data down into row 4.
$summary.Range("B4").Value = $summary.Range("B3").Value2
$summary.Range("C4").Value = $summary.Range("C3").Value2
$summary.Range("D4").Value = $summary.Range("D3").Value2

# Push row 2 (2022-Q1) data down into row 3.
$summary.Range("B3").Value = $summary.Range("B2").Value2
$summary.Range("C3").Value = $summary.Range("C2").Value2
$summary.Range("D3").Value = $summary.Range("D2").Value2

# New row 2 holds the 2022-Q3 figures.
$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 4
$summary.Range("D2").Value = 0.03

# Renumber the leading index column (0,1,2,3).
$summary.Range("A2").Value = 0
$summary.Range("A3").Value = 1
$summary.Range("A4").Value = 2
$summary.Range("A5").Value = 3

# ---------------------------------------------------------------------------
# 2) Insert a brand new "2022-Q3" worksheet right after "总计" (before the
#    existing "2022-Q1" tab) and fill it with the fund holdings table.
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Add($wb.Worksheets.Item(2))
$q3.Name = "2022-Q3"

$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

# Leading zeros / trailing zeros in these codes & percentages must survive as
# literal text, so force-text them with a quote prefix (same trick a user
# typing into Excel would use) instead of letting auto-detect coerce them to
# numbers and drop the formatting-significant digits.
$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "'000059"
$q3.Range("C2").Value = "国联安中证医药100指数A"
$q3.Range("D2").Value = "'1.70"
$q3.Range("E2").Value = "'92.19"
$q3.Range("F2").Value = "'1.12"
$q3.Range("G2").Value = "'0.0190"
$q3.Range("H2").Value = 10

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "'001351"
$q3.Range("C3").Value = "诺安中证500指数增强A"
$q3.Range("D3").Value = "'0.39"
$q3.Range("E3").Value = "'94.11"
$q3.Range("F3").Value = "'1.91"
$q3.Range("G3").Value = "'0.0074"
$q3.Range("H3").Value = 6

$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "'006569"
$q3.Range("C4").Value = "国联安中证医药100指数C"
$q3.Range("D4").Value = "'0.32"
$q3.Range("E4").Value = "'92.19"
$q3.Range("F4").Value = "'1.12"
$q3.Range("G4").Value = "'0.0036"
$q3.Range("H4").Value = 10

$q3.Range("A5").Value = 3
$q3.Range("B5").Value = "'010355"
$q3.Range("C5").Value = "诺安中证500指数增强C"
$q3.Range("D5").Value = "'0.16"
$q3.Range("E5").Value = "'94.11"
$q3.Range("F5").Value = "'1.91"
$q3.Range("G5").Value = "'0.0031"
$q3.Range("H5").Value = 6

# Match the header / leading-column styling used on the other quarter sheets.
$src = $wb.Worksheets.Item(3)
$src.Range("B1:H1").Copy()
$q3.Range("B1:H1").PasteSpecial(-4122)
$src.Range("A2").Copy()
$q3.Range("A2:A5").PasteSpecial(-4122)

$q3.Range("A1").Select()

Write-Output "done"
